# Generate Report for Handoff
#
# A new handoff report was generated. All files that are still pending
# (status "Ready for handoff") or are stuck in a failed handback
# ("Handback transform failed") get a fresh "Latest Handoff" timestamp
# recorded for every locale, because they were (re-)queued for handoff
# as part of this report run.
#
#   - Overview sheet:        column D "Latest Handoff Date"  -> 2016-03-24 22:29:40
#   - zh-cn sheet:            column E "Latest Handoff Datetime" -> 2016-03-24 22:29:35
#   - de-de sheet:            column E "Latest Handoff Datetime" -> 2016-03-24 22:29:40
#
# Affected rows (by Source File Name) are the same across all three sheets:
#   7  -> b1e5a64d-a59f-4cf4-aaf1-79f6c9e7b95c.md   (Handback transform failed)
#   10 -> 1f69bce0-4402-4edc-afff-b8355aaa931c.md   (Ready for handoff)
#   11 -> 35124e75-278f-4da3-85f4-9152a973428a.md   (Ready for handoff)
#   12 -> 558d2148-efc8-40e2-a2ee-97ce739f38fe.md   (Ready for handoff)
#   13 -> 5fe52422-bcc3-4e0c-a8e5-070574fc2395.md   (Ready for handoff)
#   14 -> 95b95732-d9ee-47dd-bbc9-4d6cce58a8ae.md   (Ready for handoff)
#   15 -> c689ce7a-c705-4072-a3e4-b8073c5398af.md   (Ready for handoff)
#   16 -> e091c753-145f-478d-ad64-0d8f4e293215.md   (Ready for handoff)

$wb = $excel.ActiveWorkbook

$affectedRows = @(7, 10, 11, 12, 13, 14, 15, 16)

# Overview sheet - "Latest Handoff Date" lives in column D
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $affectedRows) {
    $wsOverview.Range("D$r").Value = "2016-03-24 22:29:40"
}

# zh-cn sheet - "Latest Handoff Datetime" lives in column E
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $affectedRows) {
    $wsZhCn.Range("E$r").Value = "2016-03-24 22:29:35"
}

# de-de sheet - "Latest Handoff Datetime" lives in column E
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $affectedRows) {
    $wsDeDe.Range("E$r").Value = "2016-03-24 22:29:40"
}
